$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New pinout rows for a verified circuit board (HX711 load cell amp,
# --- optical counter, accelerometer) inserted above the existing
# --- "Relay Pin" / H-Bridge table ---

# Row 40: HX711 CLK -> D_IN_2
$ws.Range("C40").Value = "CLK"
$ws.Range("D40").Value = "HX711"
$ws.Range("E40").Value = "D_IN_2"

# Row 41: HX711 DOUT -> D_IN_5
$ws.Range("C41").Value = "DOUT"
$ws.Range("D41").Value = "HX711"
$ws.Range("E41").Value = "D_IN_5"

# Row 42: +5V -> Arduino (quote-prefixed like the other "+5V"/"+12V" cells)
$ws.Range("C42").Value = "'+5V"
$ws.Range("D42").Value = "Arduino"

# Row 43: Optical Counter -> D_IN_3
$ws.Range("C43").Value = "Optical Counter"
$ws.Range("D43").Value = "Optical Counter"
$ws.Range("E43").Value = "D_IN_3"

# Row 44: Accelerometer SDA -> D_IN_20
$ws.Range("C44").Value = "SDA"
$ws.Range("D44").Value = "Accelerometer"
$ws.Range("E44").Value = "D_IN_20"

# Row 45: Accelerometer SCL -> D_IN_21
$ws.Range("C45").Value = "SCL"
$ws.Range("D45").Value = "Accelerometer"
$ws.Range("E45").Value = "D_IN_21"

# --- Pin reassignment in the relay/H-bridge controller table ---
# "Start stop" Arduino pin changed from D_OUT_49 to D_OUT_11
$ws.Range("E47").Value = "D_OUT_11"

# --- Scroll / selection state as left by the author at end of session ---
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$ws.Range("H41").Select()
